$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.260865688323975
$ws.Range("B1").Value = 2.242712259292603
$ws.Range("C1").Value = 6.203773021697998
$ws.Range("D1").Value = 1.444786190986633
$ws.Range("E1").Value = 1.350248336791992
